$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the header formatting already used by C1:E1 (bold, centered, thin border)
$ws.Range("C1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill the new columns with boolean FALSE for each data row
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}
